$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '68.555.04'
Set-TextValue 'E2' '  -0.86%  '
Set-TextValue 'D3' '3.902.85'
Set-TextValue 'E3' '  +2.34%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '602.62'
Set-TextValue 'E5' '  +0.13%  '
Set-TextValue 'D6' '166.99'
Set-TextValue 'E6' '  +1.75%  '
Set-TextValue 'D7' '3.901.86'
Set-TextValue 'E7' '  +2.38%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'E9' '  -1.28%  '
Set-TextValue 'E10' '  -1.47%  '
Set-TextValue 'E11' '  +2.29%  '
Set-TextValue 'E12' '  -0.19%  '
Set-TextValue 'E13' '  +3.72%  '
Set-TextValue 'D14' '37.45'
Set-TextValue 'E14' '  +0.51%  '
Set-TextValue 'D15' '4.557.04'
Set-TextValue 'E15' '  +2.36%  '
Set-TextValue 'D16' '3.899.07'
Set-TextValue 'E16' '  +2.26%  '
Set-TextValue 'D17' '68.689.85'
Set-TextValue 'E17' '  -0.84%  '
Set-TextValue 'D18' '7.47'
Set-TextValue 'E18' '  +0.37%  '
Set-TextValue 'D19' '17.32'
Set-TextValue 'E19' '  -0.40%  '
Set-TextValue 'E20' '  -2.31%  '
Set-TextValue 'D21' '11.05'
Set-TextValue 'E21' '  -4.19%  '
Set-TextValue 'D22' '488.65'
Set-TextValue 'E23' '  +0.38%  '
Set-TextValue 'D24' '0.0000165'
Set-TextValue 'E24' '  +5.29%  '
Set-TextValue 'D25' '84.67'
Set-TextValue 'E25' '  -0.28%  '
Set-TextValue 'E26' '  -1.05%  '
Set-TextValue 'D27' '12.04'
Set-TextValue 'E27' '  -1.76%  '
Set-TextValue 'E28' '  +1.06%  '
Set-TextValue 'E29' '  +0.00%  '
Set-TextValue 'E30' '  -1.19%  '
Set-TextValue 'D31' '4.055.16'
Set-TextValue 'E31' '  +2.37%  '
Set-TextValue 'E32' '  -0.98%  '
Set-TextValue 'D33' '7.73'
Set-TextValue 'E33' '  -3.71%  '
Set-TextValue 'D34' '31.85'
Set-TextValue 'E34' '  -0.40%  '
Set-TextValue 'D35' '3.857.92'
Set-TextValue 'E35' '  +2.58%  '
Set-TextValue 'E36' '  -0.45%  '
Set-TextValue 'E37' '  +1.70%  '
Set-TextValue 'D38' '5.94'
Set-TextValue 'E38' '  +0.30%  '
Set-TextValue 'E39' '  -0.97%  '
Set-TextValue 'D40' '3.20'
Set-TextValue 'E40' '  +5.96%  '
Set-TextValue 'E41' '  +0.05%  '
Set-TextValue 'D42' '0.317'
Set-TextValue 'E42' '  -1.34%  '
Set-TextValue 'D43' '429.55'
Set-TextValue 'E43' '  +0.88%  '
Set-TextValue 'E44' '  -0.19%  '
Set-TextValue 'D45' '48.11'
Set-TextValue 'E45' '  -1.03%  '
Set-TextValue 'D46' '8.53'
Set-TextValue 'E46' '  +1.62%  '
Set-TextValue 'D48' '142.97'
Set-TextValue 'E48' '  +1.14%  '
Set-TextValue 'D49' '2.807.06'
Set-TextValue 'E49' '  -0.92%  '
Set-TextValue 'E50' '  +0.28%  '
Set-TextValue 'D51' '0.000266'
Set-TextValue 'E51' '  +17.19%  '
